$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated header row that was accidentally re-inserted at row 106.
# Deleting the entire row shifts all subsequent rows (107-196) up by one (106-195).
$ws.Rows.Item(106).Delete()
